# "Contact page worked up, Calendar page started"
#
# The calendar sheet is rebuilt from scratch (its internal sheetId moves from
# 3 -> 6 in the target, which is what happens when the sheet is deleted and a
# fresh one is inserted in its place) and populated with a first pass at a
# calendar/event table. The previously-active "deals" tab loses focus to the
# newly active "calendar" tab.

$wb = $excel.ActiveWorkbook

$callsSheet = $wb.Worksheets.Item("calls")
$oldCalendar = $wb.Worksheets.Item("calendar")
$oldCalendar.Delete() | Out-Null

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $callsSheet)
$ws.Name = "calendar"
$ws.Activate()

# Header row + sample event row. Write order below is deliberate (matches
# the original authoring order of the shared-string table) rather than a
# plain left-to-right, row-by-row fill.
$ws.Cells.Item(1, 1).Value = "title"
$ws.Cells.Item(1, 1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(2, 2).Value = "social"
$ws.Cells.Item(2, 1).Value = "Bosco Place"
$ws.Cells.Item(2, 3).Value = "re union"
$ws.Cells.Item(2, 4).Value = "Ferndale, Michigan"
$ws.Cells.Item(1, 4).Value = "location"
$ws.Cells.Item(1, 5).Value = "interval"
$ws.Cells.Item(1, 6).Value = "days"
$ws.Cells.Item(2, 6).Value = "Tuesday"
$ws.Cells.Item(2, 5).Value = "Monthly"
$ws.Cells.Item(1, 7).Value = "addNote"
$ws.Cells.Item(2, 7).Value = "Remember Me"
$ws.Cells.Item(1, 2).Value = "category"
$ws.Cells.Item(1, 3).Value = "description"

# Column widths tuned to fit their contents
$ws.Columns.Item(3).ColumnWidth = 10.09
$ws.Columns.Item(4).ColumnWidth = 17.25
$ws.Columns.Item(7).ColumnWidth = 16.76

$ws.Columns.Item(7).Select() | Out-Null
